$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 505.9655
$ws.Range("I15").Value = 505.9655
$ws.Range("K15").Value = 1517.8965
$ws.Range("M15").Value = -1348.8965
$ws.Range("H52").Value = 28031.166
$ws.Range("I52").Value = 125106.25
$ws.Range("J52").Value = 295.42856
$ws.Range("K52").Value = 375318.75
$ws.Range("L52").Value = 886.28568
$ws.Range("M52").Value = -375158.75
$ws.Range("N52").Value = -1206.28568
$ws.Range("H92").Value = 1647.6923
$ws.Range("I92").Value = 553.4
$ws.Range("J92").Value = 2331.625
$ws.Range("K92").Value = 553.4
$ws.Range("L92").Value = 2331.625
$ws.Range("M92").Value = 694.6
$ws.Range("N92").Value = -4827.625
$ws.Range("H96").Value = 2140
$ws.Range("I96").Value = 1960
$ws.Range("J96").Value = 2500
$ws.Range("K96").Value = 5880
$ws.Range("L96").Value = 7500
$ws.Range("M96").Value = -4507
$ws.Range("N96").Value = -10246
$ws.Range("H99").Value = 185
$ws.Range("I99").Value = 189
$ws.Range("J99").Value = 177
$ws.Range("K99").Value = 567
$ws.Range("L99").Value = 531
$ws.Range("M99").Value = 931
$ws.Range("N99").Value = -3527
$ws.Range("H100").Value = 5164
$ws.Range("I100").Value = 1749.8667
$ws.Range("K100").Value = 1749.8667
$ws.Range("M100").Value = -1208.8667
$ws.Range("H101").Value = 587.2
$ws.Range("I101").Value = 484.25
$ws.Range("K101").Value = 1452.75
$ws.Range("M101").Value = 169.25
$ws.Range("H103").Value = 739.7143
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 739.7143
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 2219.1429
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -3391.1429
$ws.Range("H104").Value = 2144.3333
$ws.Range("I104").Value = 2144.3333
$ws.Range("K104").Value = 6432.999899999999
$ws.Range("M104").Value = -4685.999899999999
$ws.Range("H106").Value = 2387.1738
$ws.Range("I106").Value = 2500.5625
$ws.Range("J106").Value = 2128
$ws.Range("K106").Value = 2500.5625
$ws.Range("L106").Value = 2128
$ws.Range("M106").Value = -1869.5625
$ws.Range("N106").Value = -3390
$ws.Range("H132").Value = 2154.375
$ws.Range("I132").Value = 1375.25
$ws.Range("K132").Value = 4125.75
$ws.Range("M132").Value = -1595.75
$ws.Range("H133").Value = 59514.31
$ws.Range("J133").Value = 59514.31
$ws.Range("L133").Value = 59514.31
$ws.Range("N133").Value = -69634.31
$ws.Range("H138").Value = 5367.346
$ws.Range("I138").Value = 4091.6155
$ws.Range("K138").Value = 12274.8465
$ws.Range("M138").Value = -7134.8465
$ws.Range("H141").Value = 3313.7334
$ws.Range("I141").Value = 1571.1
$ws.Range("K141").Value = 4713.299999999999
$ws.Range("M141").Value = 466.7000000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 59400
$ws.Range("J37").Value = 59400
$ws.Range("L37").Value = 59400
$ws.Range("N37").Value = -59946
$ws.Range("H44").Value = 96000
$ws.Range("J44").Value = 96000
$ws.Range("L44").Value = 96000
$ws.Range("N44").Value = -96976
$ws.Range("H55").Value = 48000
$ws.Range("J55").Value = 62000
$ws.Range("L55").Value = 62000
$ws.Range("N55").Value = -62630
$ws.Range("H97").Value = 1394.3214
$ws.Range("I97").Value = 1628.3334
$ws.Range("J97").Value = 692.2857
$ws.Range("K97").Value = 1628.3334
$ws.Range("L97").Value = 692.2857
$ws.Range("M97").Value = -1132.3334
$ws.Range("N97").Value = -1684.2857
$ws.Range("H102").Value = 3178.6924
$ws.Range("I102").Value = 3026.9167
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 3026.9167
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -1404.9167
$ws.Range("N102").Value = -8244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H55").Value = 77046.75
$ws.Range("J55").Value = 77046.75
$ws.Range("L55").Value = 77046.75
$ws.Range("N55").Value = -77592.75
$ws.Range("H94").Value = 431.5
$ws.Range("I94").Value = 336.27585
$ws.Range("K94").Value = 336.27585
$ws.Range("M94").Value = 114.72415
$ws.Range("H99").Value = 3328
$ws.Range("I99").Value = 3193.1765
$ws.Range("J99").Value = 4474
$ws.Range("K99").Value = 3193.1765
$ws.Range("L99").Value = 4474
$ws.Range("M99").Value = -1695.1765
$ws.Range("N99").Value = -7470
$ws.Range("H105").Value = 3698.5293
$ws.Range("I105").Value = 1470.2
$ws.Range("J105").Value = 4627
$ws.Range("K105").Value = 1470.2
$ws.Range("L105").Value = 4627
$ws.Range("M105").Value = 276.8
$ws.Range("N105").Value = -8121

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 445.07144
$ws.Range("I22").Value = 410.07693
$ws.Range("K22").Value = 410.07693
$ws.Range("M22").Value = -60.07693
$ws.Range("H37").Value = 14264.25
$ws.Range("J37").Value = 14264.25
$ws.Range("L37").Value = 14264.25
$ws.Range("N37").Value = -14478.25
$ws.Range("H55").Value = 18027
$ws.Range("I55").Value = 18500
$ws.Range("K55").Value = 18500
$ws.Range("M55").Value = -18185
$ws.Range("H105").Value = 1742.1
$ws.Range("I105").Value = 1742.1
$ws.Range("K105").Value = 1742.1
$ws.Range("M105").Value = 4.900000000000091
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6968.154
$ws.Range("I56").Value = 6968.154
$ws.Range("K56").Value = 6968.154
$ws.Range("M56").Value = -6438.154
$ws.Range("H119").Value = 1555.4286
$ws.Range("I119").Value = 917.8
$ws.Range("K119").Value = 2753.4
$ws.Range("M119").Value = 2084.6
$ws.Range("H133").Value = 8839.666999999999
$ws.Range("I133").Value = 12009.5
$ws.Range("J133").Value = 2500
$ws.Range("K133").Value = 36028.5
$ws.Range("L133").Value = 7500
$ws.Range("M133").Value = -30968.5
$ws.Range("N133").Value = -17620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 41000.2
$ws.Range("J69").Value = 41000.2
$ws.Range("L69").Value = 41000.2
$ws.Range("N69").Value = -42498.2
$ws.Range("H72").Value = 41000.2
$ws.Range("J72").Value = 41000.2
$ws.Range("L72").Value = 123000.6
$ws.Range("N72").Value = -130488.6
$ws.Range("H97").Value = 5437.8
$ws.Range("I97").Value = 7602.2144
$ws.Range("K97").Value = 7602.2144
$ws.Range("M97").Value = -7106.2144
$ws.Range("H99").Value = 4735.3335
$ws.Range("I99").Value = 4735.3335
$ws.Range("K99").Value = 4735.3335
$ws.Range("M99").Value = -2489.3335
$ws.Range("H113").Value = 794986.7
$ws.Range("I113").Value = 1020483.2
$ws.Range("K113").Value = 1020483.2
$ws.Range("M113").Value = -1018313.2
$ws.Range("H122").Value = 5666.7827
$ws.Range("I122").Value = 5522.467
$ws.Range("J122").Value = 5937.375
$ws.Range("K122").Value = 16567.401
$ws.Range("L122").Value = 17812.125
$ws.Range("M122").Value = -14117.401
$ws.Range("N122").Value = -22712.125
$ws.Range("H126").Value = 90912150
$ws.Range("J126").Value = 2931.6667
$ws.Range("L126").Value = 8795.000100000001
$ws.Range("N126").Value = -13735.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 563.8570999999999
$ws.Range("I16").Value = 563.8570999999999
$ws.Range("K16").Value = 563.8570999999999
$ws.Range("M16").Value = -393.8570999999999
$ws.Range("H44").Value = 5023.5
$ws.Range("J44").Value = 5023.5
$ws.Range("L44").Value = 5023.5
$ws.Range("N44").Value = -5935.5
$ws.Range("H93").Value = 1217.8182
$ws.Range("I93").Value = 1229.6
$ws.Range("K93").Value = 1229.6
$ws.Range("M93").Value = 18.40000000000009
$ws.Range("H100").Value = 11299.889
$ws.Range("J100").Value = 37000
$ws.Range("L100").Value = 37000
$ws.Range("N100").Value = -38082
$ws.Range("H132").Value = 3145.725
$ws.Range("I132").Value = 2244.7188
$ws.Range("K132").Value = 6734.1564
$ws.Range("M132").Value = -4204.1564

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 85273.086
$ws.Range("I96").Value = 112919.78
$ws.Range("K96").Value = 112919.78
$ws.Range("M96").Value = -111546.78
$ws.Range("H121").Value = 90210
$ws.Range("J121").Value = 90210
$ws.Range("L121").Value = 90210
$ws.Range("N121").Value = -93704
$ws.Range("H132").Value = 3702.6038
$ws.Range("I132").Value = 2893.9211
$ws.Range("K132").Value = 8681.763300000001
$ws.Range("M132").Value = -6151.763300000001
